$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.350.56'
$ws.Range("E2").Value = '  +2.37%  '
$ws.Range("D3").Value = '2.694.52'
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.73'
$ws.Range("E5").Value = '  +1.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.96'
$ws.Range("E6").Value = '  +1.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.577'
$ws.Range("E8").Value = '  +1.70%  '
$ws.Range("D9").Value = '2.713.22'
$ws.Range("E9").Value = '  +0.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.52'
$ws.Range("E10").Value = '  +4.10%  '
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("E12").Value = '  +1.51%  '
$ws.Range("E13").Value = '  +1.63%  '
$ws.Range("D14").Value = '3.170.93'
$ws.Range("E14").Value = '  +1.28%  '
$ws.Range("D15").Value = '60.365.35'
$ws.Range("E15").Value = '  +2.36%  '
$ws.Range("D16").Value = '2.852.75'
$ws.Range("E16").Value = '  +6.38%  '
$ws.Range("E17").Value = '  +1.84%  '
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '352.26'
$ws.Range("E19").Value = '  +1.04%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.60'
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.36'
$ws.Range("E22").Value = '  +3.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.17'
$ws.Range("E24").Value = '  +3.56%  '
$ws.Range("E25").Value = '  +0.81%  '
$ws.Range("E26").Value = '  +5.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.37'
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0821'
$ws.Range("E29").Value = '  +1.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.91'
$ws.Range("E30").Value = '  +7.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("E32").Value = '  +1.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.17'
$ws.Range("E33").Value = '  +0.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '147.71'
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.32'
$ws.Range("E35").Value = '  +6.76%  '
$ws.Range("E36").Value = '  +8.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.953'
$ws.Range("E37").Value = '  -5.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.53'
$ws.Range("E38").Value = '  +9.14%  '
$ws.Range("E39").Value = '  +3.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.95'
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '287.00'
$ws.Range("E42").Value = '  +3.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.17'
$ws.Range("E43").Value = '  +2.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.616'
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.146.23'
$ws.Range("E46").Value = '  +6.85%  '
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.995'
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.90'
$ws.Range("E48").Value = '  +3.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0539'
$ws.Range("E49").Value = '  +0.89%  '
$ws.Range("E50").Value = '  +1.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.45'
$ws.Range("E51").Value = '  +1.69%  '
